$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting existing rows 119:185 down to 120:186
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new data record
$ws.Cells.Item(119,1).Value2 = 1
$ws.Cells.Item(119,2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(119,3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(119,4).Value2 = 44455
$ws.Cells.Item(119,5).Value2 = 15
$ws.Cells.Item(119,6).Value2 = 100114013
$ws.Cells.Item(119,7).Value2 = "Zanahoria"
$ws.Cells.Item(119,8).Value2 = "Sin especificar"
$ws.Cells.Item(119,9).Value2 = "Primera"
$ws.Cells.Item(119,10).Value2 = 60
$ws.Cells.Item(119,11).Value2 = 7500
$ws.Cells.Item(119,12).Value2 = 8000
$ws.Cells.Item(119,13).Value2 = 7750
$ws.Cells.Item(119,14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(119,15).Value2 = "Valle de Camiña"
$ws.Cells.Item(119,16).Value2 = 310
$ws.Cells.Item(119,17).Value2 = 25
$ws.Cells.Item(119,18).Value2 = "Hortaliza"
